# Generate Report for Handoff
# Re-running the handoff-status report bumped a couple of timestamps and
# filled in the default "ht" priority for rows that previously had no
# explicit priority recorded ("low").

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 -> Priority (E) "low" -> "ht",
#              Latest Handoff Datetime (H) bumped 30s
foreach ($row in 4..7) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-09-01 20:36:19"
}

# de-de sheet: rows 4-7 -> Priority (E) "low" -> "ht"
foreach ($row in 4..7) {
    $wsDeDe.Range("E$row").Value = "ht"
}

# Shared "Latest HO Xliff Generate Date"/"Latest Handoff Datetime" timestamp
# bumped 31s -- appears on Overview!G4:G7 and de-de!H4:H7 (same shared string).
foreach ($row in 4..7) {
    $wsOverview.Range("G$row").Value = "2016-09-01 20:36:24"
    $wsDeDe.Range("H$row").Value = "2016-09-01 20:36:24"
}
